# GMS Data Release 1
# Sequencing_report data dictionary update: rework the field list to match
# the new GMS schema (platekey / referral_id / associated_interpretation_
# request_id / delivery_type / data_format) and refresh the view/print setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "lab_sample_id" row (row 3) becomes the "platekey" row ---
$ws.Range("B3").Value = "platekey"
$ws.Range("C3").Value = "varchar"
$ws.Range("D3").Value = "Concatination of Plate ID and Well ID - unique identifier for a proccessed well"

# --- 2. Insert three new rows after row 3 for the new fields ---
$ws.Rows.Item(4).Resize(3, 1).Insert()

$ws.Range("A4").Value = "sequencing_report"
$ws.Range("B4").Value = "referral_id"
$ws.Range("C4").Value = "varchar"
$ws.Range("D4").Value = "Varchar"

$ws.Range("A5").Value = "sequencing_report"
$ws.Range("B5").Value = "associated_interpretation_request_id"
$ws.Range("C5").Value = "varchar"
$ws.Range("D5").Value = "Varchar"

$ws.Range("A6").Value = "sequencing_report"
$ws.Range("B6").Value = "delivery_type"
$ws.Range("C6").Value = "[rare disease germline, cancer germline, cancer somatic]"
$ws.Range("D6").Value = "Type the sample type "

# match the table's existing body-cell formatting (Arial 12, wrap, thin
# border, left/top aligned) on the freshly-inserted rows
$newRng = $ws.Range("A4:D6")
$newRng.Font.Name = "Arial"
$newRng.Font.Size = 12
$newRng.HorizontalAlignment = -4131
$newRng.VerticalAlignment = -4160
$newRng.WrapText = $true
$newRng.Borders.LineStyle = 1
$newRng.Borders.Weight = 2

# --- 3. Drop the old "plate_key" and "type" rows - superseded above ---
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# --- 4. Append a new "data_format" row at the bottom of the table ---
$ws.Range("A12").Value = "sequencing_report"
$ws.Range("B12").Value = "data_format"
$ws.Range("C12").Value = "varchar"
$ws.Range("D12").Value = ""

$lastRng = $ws.Range("A12:D12")
$lastRng.Font.Name = "Arial"
$lastRng.Font.Size = 12
$lastRng.HorizontalAlignment = -4131
$lastRng.VerticalAlignment = -4160
$lastRng.WrapText = $true
$lastRng.Borders.LineStyle = 1
$lastRng.Borders.Weight = 2

# --- 5. Row heights (auto-computed by Excel on open; set explicitly here) ---
$ws.Rows.Item(1).RowHeight = 31
$ws.Rows.Item(2).RowHeight = 155
$ws.Rows.Item(3).RowHeight = 170.5
$ws.Rows.Item(4).RowHeight = 46.5
$ws.Rows.Item(5).RowHeight = 46.5
$ws.Rows.Item(6).RowHeight = 62
$ws.Rows.Item(7).RowHeight = 155
$ws.Rows.Item(8).RowHeight = 108.5
$ws.Rows.Item(9).RowHeight = 77.5
$ws.Rows.Item(10).RowHeight = 201.5
$ws.Rows.Item(11).RowHeight = 232.5
$ws.Rows.Item(12).RowHeight = 46.5

# --- 6. Column widths tweaked slightly; new narrow column E ---
$ws.Columns.Item(2).ColumnWidth = 22.26953125
$ws.Columns.Item(3).ColumnWidth = 36.1796875
$ws.Columns.Item(5).ColumnWidth = 8.7265625

# --- 7. Default row height / font metrics for the new Excel version ---
$ws.Application.StandardHeight = 14.5

# --- 8. View state: scrolled down a couple of rows, new selection ---
$ws.Range("G6").Select()

# --- 9. Page setup now defined for printing ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
